# Anpassung Technische Risiken gemaess Review
# - neue Spalte I "Kommentar" mit Review-Kommentaren je Risiko
# - Zeilenhoehen fuer Zeilen 8, 9, 13 angepasst
# - Seiteneinrichtung: FitToPage + Skalierung 71%
# - Fenster-/Auswahlzustand aktualisiert

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- neue Spalte I: Breite ---
$ws.Columns.Item(9).ColumnWidth = 18.83

# --- Kopfzeile Spalte I ---
$ws.Range("I7").Value = "Kommentar"

# --- Kommentare je Risiko (Spalte I, Zeilen 8-13) ---
$comments = @{
    8  = "Durch Vereinfachung des Konzept (z.B.: RabbitMQ einsetzen) soweit umgangen."
    9  = "Wurde getstet durch ein Menge Messages an einen RabbitMQ Server zu senden -> funtkioniert"
    10 = "Performance erfüllt unsere Anforderungen"
    11 = "Alles soweit vereinfacht wie möglich -> funktioniert"
    12 = "Backup kann auf einem anderen Server eingespielt werden."
    13 = "Bereits schon mehrmals passiert durch wiederherstellen aus alten Commits kein Problem."
}

foreach ($row in $comments.Keys) {
    $cell = $ws.Range("I$row")
    $cell.Value = $comments[$row]
    $cell.WrapText = $true
    $cell.VerticalAlignment = -4160
}

# --- Zeilenhoehen (Autofit-Ergebnis der neuen, laengeren Kommentare) ---
$ws.Rows.Item(8).RowHeight = 60
$ws.Rows.Item(9).RowHeight = 75
$ws.Rows.Item(13).RowHeight = 90

# --- Seiteneinrichtung ---
$ws.EnableFormatConditionsCalculation = $false
$ws.PageSetup.Zoom = 71
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1

# --- Fensterposition / Bildlauf / Auswahl ---
$wb.Windows.Item(1).Left = 120
$wb.Windows.Item(1).Top = 0
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I13").Select() | Out-Null
